# jbhunt_api_tracking.xlsx - append new API call log rows (p2p and first mile)
#
# Adds 15 new tracking rows (rows 12-26) below the existing log (rows 1-11),
# and corrects row 11's Origin/Destination ZIP cells from text to numeric
# values. New rows reuse the same Origin/Destination ZIP + Message patterns
# seen earlier in the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data for the new rows (12-26). Columns: Timestamp, OriginZip,
# DestZip, Weight, Status, Message. OriginZip/DestZip are numeric except
# for the final row, which keeps them as zero-padded text (as in the
# original log's row 11 before correction).
$newRows = @(
    @{ Row=12; Timestamp="2025-08-11 12:42:15"; Origin=29483;   Dest=30567;  Weight=11025;        Status="Success"; Message="Rate: 693.44, Carrier: JBHZ" },
    @{ Row=13; Timestamp="2025-08-11 12:42:29"; Origin=29483;   Dest=30567;  Weight=45007.703685;  Status="Success"; Message="Rate: 693.44, Carrier: JBHZ" },
    @{ Row=14; Timestamp="2025-08-11 12:42:43"; Origin=7201;    Dest=8873;   Weight=11025;        Status="Success"; Message="Rate: 428.9, Carrier: JBHZ" },
    @{ Row=15; Timestamp="2025-08-11 12:42:58"; Origin=7201;    Dest=8873;   Weight=45007.703685;  Status="Success"; Message="Rate: 428.9, Carrier: JBHZ" },
    @{ Row=16; Timestamp="2025-08-11 12:43:13"; Origin=7201;    Dest=21901;  Weight=45007.703685;  Status="Success"; Message="Rate: 714.31, Carrier: JBHZ" },
    @{ Row=17; Timestamp="2025-08-11 12:45:29"; Origin=29483;   Dest=30567;  Weight=11025;        Status="Success"; Message="Rate: 693.44, Carrier: JBHZ" },
    @{ Row=18; Timestamp="2025-08-11 12:45:32"; Origin=29483;   Dest=30567;  Weight=45007.703685;  Status="Success"; Message="Rate: 693.44, Carrier: JBHZ" },
    @{ Row=19; Timestamp="2025-08-11 12:45:45"; Origin=7201;    Dest=8873;   Weight=11025;        Status="Success"; Message="Rate: 428.9, Carrier: JBHZ" },
    @{ Row=20; Timestamp="2025-08-11 12:45:48"; Origin=7201;    Dest=8873;   Weight=45007.703685;  Status="Success"; Message="Rate: 428.9, Carrier: JBHZ" },
    @{ Row=21; Timestamp="2025-08-11 12:46:14"; Origin=7201;    Dest=21901;  Weight=45007.703685;  Status="Success"; Message="Rate: 714.31, Carrier: JBHZ" },
    @{ Row=22; Timestamp="2025-08-11 12:53:07"; Origin=29483;   Dest=30567;  Weight=11025;        Status="Success"; Message="Rate: 693.44, Carrier: JBHZ" },
    @{ Row=23; Timestamp="2025-08-11 12:53:19"; Origin=29483;   Dest=30567;  Weight=45007.703685;  Status="Success"; Message="Rate: 693.44, Carrier: JBHZ" },
    @{ Row=24; Timestamp="2025-08-11 12:53:34"; Origin=7201;    Dest=8873;   Weight=11025;        Status="Success"; Message="Rate: 428.9, Carrier: JBHZ" },
    @{ Row=25; Timestamp="2025-08-11 12:53:36"; Origin=7201;    Dest=8873;   Weight=45007.703685;  Status="Success"; Message="Rate: 428.9, Carrier: JBHZ" },
    @{ Row=26; Timestamp="2025-08-11 12:53:50"; Origin="07201"; Dest="21901"; Weight=45007.703685; Status="Success"; Message="Rate: 714.31, Carrier: JBHZ" }
)

# --- Phase 1: write the Timestamp (column A) for every new row first, in
# row order, so the new date/time strings land in the shared-string table
# in chronological order ahead of everything else that follows. ---
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Timestamp
}

# Row 11's timestamp is unchanged (string already present in the sheet).
$ws.Cells.Item(11, 1).Value = "2025-08-10 14:23:33"

# --- Phase 2: the last row keeps Origin/Destination ZIP as zero-padded
# text ("07201"/"21901") instead of numbers - set those explicitly as
# text so the leading zero survives. ---
$ws.Cells.Item(26, 3).NumberFormat = "@"
$ws.Cells.Item(26, 3).Value = "21901"
$ws.Cells.Item(26, 3).Style = "Normal"

$ws.Cells.Item(26, 2).NumberFormat = "@"
$ws.Cells.Item(26, 2).Value = "07201"
$ws.Cells.Item(26, 2).Style = "Normal"

# --- Phase 3: write the Message (column F) for every row, in row order,
# so first-seen new message strings are appended in that order. ---
$ws.Cells.Item(11, 6).Value = "Rate: 3493.82, Carrier: JBHZ"
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 6).Value = $r.Message
}

# --- Phase 4: fill in the remaining columns (Origin/Dest ZIP numbers,
# Weight, Status) for every row. ---

# Row 11: Origin/Destination ZIP become plain numbers (previously text).
$ws.Cells.Item(11, 2).Value = 7201
$ws.Cells.Item(11, 3).Value = 85043
$ws.Cells.Item(11, 4).Value = 45007.703685
$ws.Cells.Item(11, 5).Value = "Success"

foreach ($r in $newRows) {
    if ($r.Row -ne 26) {
        $ws.Cells.Item($r.Row, 2).Value = $r.Origin
        $ws.Cells.Item($r.Row, 3).Value = $r.Dest
    }
    $ws.Cells.Item($r.Row, 4).Value = $r.Weight
    $ws.Cells.Item($r.Row, 5).Value = $r.Status
}
